$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf1"
$ws.Cells.Item(2,3).Value = "Cd44"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.161357666666667
$ws.Cells.Item(2,8).Value = 3.484073
$ws.Cells.Item(2,9).Value = 0.1270850363824361
$ws.Cells.Item(2,10).Value = 0.1270850363824361
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 239.0839323333333
$ws.Cells.Item(2,14).Value = 717.251797
$ws.Cells.Item(2,15).Value = 0.4086975387666237
$ws.Cells.Item(2,16).Value = 0.4086975387666237
$ws.Cells.Item(2,17).Value = 277.6619577921313
$ws.Cells.Item(2,18).Value = 2498.957620129181
$ws.Cells.Item(2,19).Value = 0.05193934158356846
$ws.Cells.Item(2,20).Value = 0.05193934158356846

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf1"
$ws.Cells.Item(3,3).Value = "Cd44"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.161357666666667
$ws.Cells.Item(3,8).Value = 3.484073
$ws.Cells.Item(3,9).Value = 0.1270850363824361
$ws.Cells.Item(3,10).Value = 0.1270850363824361
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 117.0512696666667
$ws.Cells.Item(3,14).Value = 351.153809
$ws.Cells.Item(3,15).Value = 0.2000910950200451
$ws.Cells.Item(3,16).Value = 0.2000910950200451
$ws.Cells.Item(3,17).Value = 135.9383894204508
$ws.Cells.Item(3,18).Value = 1223.445504784057
$ws.Cells.Item(3,19).Value = 0.02542858409042391
$ws.Cells.Item(3,20).Value = 0.02542858409042391

# Row 4: ECs -> M2
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf1"
$ws.Cells.Item(4,3).Value = "Cd44"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.161357666666667
$ws.Cells.Item(4,8).Value = 3.484073
$ws.Cells.Item(4,9).Value = 0.1270850363824361
$ws.Cells.Item(4,10).Value = 0.1270850363824361
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 171.15883
$ws.Cells.Item(4,14).Value = 513.47649
$ws.Cells.Item(4,15).Value = 0.2925842480357353
$ws.Cells.Item(4,16).Value = 0.2925842480357353
$ws.Cells.Item(4,17).Value = 198.7766194381967
$ws.Cells.Item(4,18).Value = 1788.98957494377
$ws.Cells.Item(4,19).Value = 0.03718307980654913
$ws.Cells.Item(4,20).Value = 0.03718307980654913

# Row 5: ECs -> sCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Fgf1"
$ws.Cells.Item(5,3).Value = "Cd44"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 1.161357666666667
$ws.Cells.Item(5,8).Value = 3.484073
$ws.Cells.Item(5,9).Value = 0.1270850363824361
$ws.Cells.Item(5,10).Value = 0.1270850363824361
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 57.695868
$ws.Cells.Item(5,14).Value = 173.087604
$ws.Cells.Item(5,15).Value = 0.09862711817759588
$ws.Cells.Item(5,16).Value = 0.09862711817759588
$ws.Cells.Item(5,17).Value = 67.005538636788
$ws.Cells.Item(5,18).Value = 603.049847731092
$ws.Cells.Item(5,19).Value = 0.0125340309018946
$ws.Cells.Item(5,20).Value = 0.0125340309018946

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf1"
$ws.Cells.Item(6,3).Value = "Cd44"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 5.168173666666667
$ws.Cells.Item(6,8).Value = 15.504521
$ws.Cells.Item(6,9).Value = 0.565542861868062
$ws.Cells.Item(6,10).Value = 0.565542861868062
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 239.0839323333333
$ws.Cells.Item(6,14).Value = 717.251797
$ws.Cells.Item(6,15).Value = 0.4086975387666237
$ws.Cells.Item(6,16).Value = 0.4086975387666237
$ws.Cells.Item(6,17).Value = 1235.627283208249
$ws.Cells.Item(6,18).Value = 11120.64554887424
$ws.Cells.Item(6,19).Value = 0.2311359757125096
$ws.Cells.Item(6,20).Value = 0.2311359757125096

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf1"
$ws.Cells.Item(7,3).Value = "Cd44"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 5.168173666666667
$ws.Cells.Item(7,8).Value = 15.504521
$ws.Cells.Item(7,9).Value = 0.565542861868062
$ws.Cells.Item(7,10).Value = 0.565542861868062
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 117.0512696666667
$ws.Cells.Item(7,14).Value = 351.153809
$ws.Cells.Item(7,15).Value = 0.2000910950200451
$ws.Cells.Item(7,16).Value = 0.2000910950200451
$ws.Cells.Item(7,17).Value = 604.9412895411655
$ws.Cells.Item(7,18).Value = 5444.47160587049
$ws.Cells.Item(7,19).Value = 0.1131600905119506
$ws.Cells.Item(7,20).Value = 0.1131600905119506

# Row 8: FAPs -> M2
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Fgf1"
$ws.Cells.Item(8,3).Value = "Cd44"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.168173666666667
$ws.Cells.Item(8,8).Value = 15.504521
$ws.Cells.Item(8,9).Value = 0.565542861868062
$ws.Cells.Item(8,10).Value = 0.565542861868062
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 171.15883
$ws.Cells.Item(8,14).Value = 513.47649
$ws.Cells.Item(8,15).Value = 0.2925842480357353
$ws.Cells.Item(8,16).Value = 0.2925842480357353
$ws.Cells.Item(8,17).Value = 884.5785580234767
$ws.Cells.Item(8,18).Value = 7961.207022211291
$ws.Cells.Item(8,19).Value = 0.1654689329716446
$ws.Cells.Item(8,20).Value = 0.1654689329716446

# Row 9: FAPs -> sCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Fgf1"
$ws.Cells.Item(9,3).Value = "Cd44"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.168173666666667
$ws.Cells.Item(9,8).Value = 15.504521
$ws.Cells.Item(9,9).Value = 0.565542861868062
$ws.Cells.Item(9,10).Value = 0.565542861868062
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 57.695868
$ws.Cells.Item(9,14).Value = 173.087604
$ws.Cells.Item(9,15).Value = 0.09862711817759588
$ws.Cells.Item(9,16).Value = 0.09862711817759588
$ws.Cells.Item(9,17).Value = 298.182265673076
$ws.Cells.Item(9,18).Value = 2683.640391057684
$ws.Cells.Item(9,19).Value = 0.05577786267195713
$ws.Cells.Item(9,20).Value = 0.05577786267195713

# Row 10: sCs -> ECs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fgf1"
$ws.Cells.Item(10,3).Value = "Cd44"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.808898333333333
$ws.Cells.Item(10,8).Value = 8.426695
$ws.Cells.Item(10,9).Value = 0.3073721017495019
$ws.Cells.Item(10,10).Value = 0.3073721017495019
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 239.0839323333333
$ws.Cells.Item(10,14).Value = 717.251797
$ws.Cells.Item(10,15).Value = 0.4086975387666237
$ws.Cells.Item(10,16).Value = 0.4086975387666237
$ws.Cells.Item(10,17).Value = 671.5624590578794
$ws.Cells.Item(10,18).Value = 6044.062131520916
$ws.Cells.Item(10,19).Value = 0.1256222214705457
$ws.Cells.Item(10,20).Value = 0.1256222214705457

# Row 11: sCs -> FAPs
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Fgf1"
$ws.Cells.Item(11,3).Value = "Cd44"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.808898333333333
$ws.Cells.Item(11,8).Value = 8.426695
$ws.Cells.Item(11,9).Value = 0.3073721017495019
$ws.Cells.Item(11,10).Value = 0.3073721017495019
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 117.0512696666667
$ws.Cells.Item(11,14).Value = 351.153809
$ws.Cells.Item(11,15).Value = 0.2000910950200451
$ws.Cells.Item(11,16).Value = 0.2000910950200451
$ws.Cells.Item(11,17).Value = 328.7851162812506
$ws.Cells.Item(11,18).Value = 2959.066046531255
$ws.Cells.Item(11,19).Value = 0.06150242041767055
$ws.Cells.Item(11,20).Value = 0.06150242041767055

# Row 12: sCs -> M2
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Fgf1"
$ws.Cells.Item(12,3).Value = "Cd44"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 2.808898333333333
$ws.Cells.Item(12,8).Value = 8.426695
$ws.Cells.Item(12,9).Value = 0.3073721017495019
$ws.Cells.Item(12,10).Value = 0.3073721017495019
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 171.15883
$ws.Cells.Item(12,14).Value = 513.47649
$ws.Cells.Item(12,15).Value = 0.2925842480357353
$ws.Cells.Item(12,16).Value = 0.2925842480357353
$ws.Cells.Item(12,17).Value = 480.7677523222833
$ws.Cells.Item(12,18).Value = 4326.90977090055
$ws.Cells.Item(12,19).Value = 0.08993223525754153
$ws.Cells.Item(12,20).Value = 0.08993223525754153

# Row 13: sCs -> sCs
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Fgf1"
$ws.Cells.Item(13,3).Value = "Cd44"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 2.808898333333333
$ws.Cells.Item(13,8).Value = 8.426695
$ws.Cells.Item(13,9).Value = 0.3073721017495019
$ws.Cells.Item(13,10).Value = 0.3073721017495019
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 57.695868
$ws.Cells.Item(13,14).Value = 173.087604
$ws.Cells.Item(13,15).Value = 0.09862711817759588
$ws.Cells.Item(13,16).Value = 0.09862711817759588
$ws.Cells.Item(13,17).Value = 162.06182746542
$ws.Cells.Item(13,18).Value = 1458.55644718878
$ws.Cells.Item(13,19).Value = 0.03031522460374415
$ws.Cells.Item(13,20).Value = 0.03031522460374415

